$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "[Aggarwal%Gaurav%coreGivesNoEmail%0,              Aggarwal%Saurabh%coreGivesNoEmail%0,              Garcia-Telles%Nelson%coreGivesNoEmail%0,              Henry%Brandon Michael%coreGivesNoEmail%0,              Lavie%Carl%coreGivesNoEmail%0,              Lippi%Giuseppe%coreGivesNoEmail%0]"
$ws.Range("E3").Value = "[Juliana%Amorim dos Santos%NULL%1,               Ana Gabriela Costa%Normando%NULL%1,               Rainier Luiz%Carvalho da Silva%NULL%1,               Renata Monteiro%De Paula%NULL%1,               Allan Christian%Cembranel%NULL%1,               Alan Roger%Santos-Silva%NULL%1,               Eliete Neves Silva%Guerra%NULL%1]"
$ws.Range("E4").Value = "[Reza%Ansari%xref no email%1,        Mina%Gheitani%xref no email%1,        Farrokh%Heidari%xref no email%1,        Firouzeh%Heidari%xref no email%1]"
$ws.Range("E5").Value = "[Á.%Beltrán‐Corbellini%beltran_corbellini@hotmail.com%2,               J. L.%Chico‐García%NULL%4,               J. L.%Chico‐García%NULL%0,               J.%Martínez‐Poles%NULL%2,               F.%Rodríguez‐Jorge%NULL%2,               E.%Natera‐Villalba%NULL%2,               J.%Gómez‐Corral%NULL%4,               J.%Gómez‐Corral%NULL%0,               A.%Gómez‐López%NULL%2,               E.%Monreal%NULL%2,               P.%Parra‐Díaz%NULL%2,               J. L.%Cortés‐Cuevas%NULL%4,               J. L.%Cortés‐Cuevas%NULL%0,               J. C.%Galán%NULL%2,               C.%Fragola‐Arnau%NULL%2,               J.%Porta‐Etessam%NULL%2,               J.%Masjuan%NULL%2,               A.%Alonso‐Cánovas%NULL%4,               A.%Alonso‐Cánovas%NULL%0]"
$ws.Range("E6").Value = "[François%Bénézit%NULL%2,               Paul%Le Turnier%NULL%2,               Charles%Declerck%NULL%2,               Cécile%Paillé%NULL%2,               Matthieu%Revest%NULL%2,               Vincent%Dubée%NULL%2,               Pierre%Tattevin%pierre.tattevin@chu-rennes.fr%2,               Cédric%Arvieux%NULL%2,               Marion%Baldeyrou%NULL%2,               Jean-Marc%Chapplain%NULL%2,               Pauline%Comacle%NULL%2,               Solène%Patrat-Delon%NULL%2,               Anne%Maillard%NULL%2,               Mélanie%Poinot%NULL%2,               Charlotte%Pronier%NULL%2,               Faouzi%Souala%NULL%2,               Vincent%Thibault%NULL%2,               Pierre%Abgueguen%NULL%2,               Hélène%Cormier%NULL%2,               Valérie%Delbos%NULL%2,               Marine%de la Chapelle%NULL%2,               Alexandra%Ducancelle%NULL%2,               Rafael%Mahieu%NULL%2,               Valérie%Rabier%NULL%2,               Sami%Rehaiem%NULL%2,               Yves%Vandamme%NULL%2,               Charlotte%Biron%NULL%2,               Jeanne%Brochon%NULL%2,               David%Boutoille%NULL%0,               Marie%Chauveau%NULL%2,               Colin%Deschanvres%NULL%2,               Benjamin J%Gaborit%NULL%2,               Joël%Jenvrin%NULL%2,               Raphaël%Lecomte%NULL%2,               Maeva%Lefebvre%NULL%2,               François%Raffi%NULL%2]"
$ws.Range("E7").Value = "[ Alex%Carignan%null%1,                Louis%Valiquette%null%1,                Cynthia%Grenier%null%1,                Jean Berchmans%Musonera%null%1,                Delphin%Nkengurutse%null%1,                Anaïs%Marcil-Héguy%null%1,                Kim%Vettese%null%1,                Dominique%Marcoux%null%1,                Corinne%Valiquette%null%1,                Wei Ting%Xiong%null%1,                Pierre-Hughes%Fortier%null%1,                Mélissa%Généreux%null%1,                Jacques%Pépin%null%1,              Alex%Carignan%null%1,              Louis%Valiquette%null%1,              Cynthia%Grenier%null%1,              Jean Berchmans%Musonera%null%1,              Delphin%Nkengurutse%null%1,              Anaïs%Marcil-Héguy%null%1,              Kim%Vettese%null%1,              Dominique%Marcoux%null%1,              Corinne%Valiquette%null%1,              Wei Ting%Xiong%null%1,              Pierre-Hughes%Fortier%null%1,              Mélissa%Généreux%null%1,              Jacques%Pépin%null%1]"
$ws.Range("E8").Value = "[Filiz%Cebeci Kahraman%cebecifiliz@yahoo.com%1,               Hülya%Çaşkurlu%NULL%2,               Hülya%Çaşkurlu%NULL%0]"
$ws.Range("E9").Value = "[Eléonore%Chary%NULL%1,               Florent%Carsuzaa%NULL%1,               Jean-Paul%Trijolet%NULL%2,               Jean-Paul%Trijolet%NULL%0,               Anne-Laure%Capitaine%NULL%1,               Mariam%Roncato-Saberan%NULL%1,               Kevin%Fouet%NULL%1,               France%Cazenave-Roblot%NULL%0,               Mélanie%Catroux%NULL%1,               Caroline%Allix-Beguec%NULL%1,               Xavier%Dufour%NULL%1]"
$ws.Range("E10").Value = "[Aline%Desoutter%coreGivesNoEmail%1,              Anne-Ga\u00eblle%Chaux-Bodard%coreGivesNoEmail%1,              Sophie%Deneuve%coreGivesNoEmail%1]"
$ws.Range("E11").Value = "[Andrea%De Maria%de-maria@unige.it%1,               Paola%Varese%NULL%2,               Paola%Varese%NULL%0,               Chiara%Dentone%NULL%1,               Emanuela%Barisione%NULL%1,               Matteo%Bassetti%NULL%1]"
$ws.Range("E12").Value = "[Valeria%Dell’Era%valeria.dellera@gmail.com%1,               Filippo%Farri%NULL%2,               Filippo%Farri%NULL%0,               Giacomo%Garzaro%NULL%1,               Miriam%Gatto%NULL%1,               Paolo%Aluffi Valletti%NULL%0,               Massimiliano%Garzaro%NULL%0]"
$ws.Range("E13").Value = "[Matteo%Gelardi%NULL%1,               Eleonora%Trecca%NULL%1,               Michele%Cassano%NULL%1,               Giorgio%Ciprandi%NULL%1]"
$ws.Range("E14").Value = "[Andrea%Giacomelli%andrea.giacomelli@unimi.it%0,               Laura%Pezzati%NULL%4,               Laura%Pezzati%NULL%0,               Federico%Conti%NULL%2,               Dario%Bernacchia%NULL%0,               Matteo%Siano%NULL%0,               Letizia%Oreni%NULL%0,               Stefano%Rusconi%NULL%0,               Cristina%Gervasoni%NULL%2,               Anna Lisa%Ridolfo%NULL%0,               Giuliano%Rizzardini%NULL%0,               Spinello%Antinori%NULL%0,               Massimo%Galli%NULL%0,               Massimo%Galli%NULL%0]"
$ws.Range("E15").Value = "[G.-u.%Kim%NULL%1,               M.-J.%Kim%NULL%1,               S.H.%Ra%NULL%1,               J.%Lee%NULL%1,               S.%Bae%NULL%1,               J.%Jung%NULL%1,               S.-H.%Kim%kimsunghanmd@hotmail.com%1]"
$ws.Range("E16").Value = "[T.%Klopfenstein%NULL%2,               N.J.%Kadiane-Oussou%NULL%2,               L.%Toko%NULL%2,               P.-Y.%Royer%NULL%2,               Q.%Lepiller%NULL%2,               V.%Gendrin%NULL%2,               S.%Zayet%NULL%2]"
$ws.Range("E17").Value = "[Marjolein F. Q.%Kluytmans-van den Bergh%NULL%1,               Anton G. M.%Buiting%NULL%1,               Suzan D.%Pas%NULL%1,               Robbert G.%Bentvelsen%NULL%1,               Wouter%van den Bijllaardt%NULL%1,               Anne J. G.%van Oudheusden%NULL%1,               Miranda M. L.%van Rijen%NULL%1,               Jaco J.%Verweij%NULL%1,               Marion P. G.%Koopmans%NULL%1,               Jan A. J. W.%Kluytmans%NULL%1]"
$ws.Range("E18").Value = "[Frédéric%Lapostolle%frederic.lapostolle@aphp.fr%1,               Elodie%Schneider%NULL%1,               Isabelle%Vianu%NULL%1,               Guillaume%Dollet%NULL%1,               Bastien%Roche%NULL%1,               Julia%Berdah%NULL%1,               Julie%Michel%NULL%1,               Laurent%Goix%NULL%1,               Erick%Chanzy%NULL%1,               Tomislav%Petrovic%NULL%1,               Frédéric%Adnet%NULL%1]"
$ws.Range("E19").Value = "[Jerome R.%Lechien%Jerome.Lechien@umons.ac.be%2,               Carlos M.%Chiesa-Estomba%NULL%2,               Daniele R.%De Siati%NULL%2,               Mihaela%Horoi%NULL%2,               Serge D.%Le Bon%NULL%2,               Alexandra%Rodriguez%NULL%2,               Didier%Dequanter%NULL%2,               Serge%Blecic%NULL%2,               Fahd%El Afia%NULL%2,               Lea%Distinguin%NULL%2,               Younes%Chekkoury-Idrissi%NULL%2,               Stéphane%Hans%NULL%2,               Irene Lopez%Delgado%NULL%2,               Christian%Calvo-Henriquez%NULL%2,               Philippe%Lavigne%NULL%2,               Chiara%Falanga%NULL%2,               Maria Rosaria%Barillari%NULL%2,               Giovanni%Cammaroto%NULL%2,               Mohamad%Khalife%NULL%2,               Pierre%Leich%NULL%2,               Christel%Souchay%NULL%2,               Camelia%Rossi%NULL%2,               Fabrice%Journe%NULL%2,               Julien%Hsieh%NULL%2,               Myriam%Edjlali%NULL%2,               Robert%Carlier%NULL%2,               Laurence%Ris%NULL%2,               Andrea%Lovato%NULL%2,               Cosimo%De Filippis%NULL%2,               Frederique%Coppee%NULL%2,               Nicolas%Fakhry%NULL%2,               Tareck%Ayad%NULL%2,               Sven%Saussez%NULL%2]"
$ws.Range("E20").Value = "[Daniel J.%Lee%NULL%1,               Joel%Lockwood%NULL%2,               Joel%Lockwood%NULL%0,               Paul%Das%NULL%1,               Ri%Wang%NULL%1,               Eitan%Grinspun%NULL%1,               John M.%Lee%NULL%1]"
$ws.Range("E21").Value = "[Yonghyun%Lee%NULL%1,               Pokkee%Min%NULL%2,               Pokkee%Min%NULL%0,               Seonggu%Lee%NULL%2,               Seonggu%Lee%NULL%0,               Shin-Woo%Kim%NULL%0,               Shin-Woo%Kim%NULL%0]"
$ws.Range("E22").Value = "[Ruth%Levinson%xref no email%1,        Meital%Elbaz%xref no email%1,        Ronen%Ben-Ami%xref no email%1,        David%Shasha%xref no email%1,        Tal%Levinson%xref no email%1,        Guy%Choshen%xref no email%1,        Ksenia%Petrov%xref no email%1,        Avi%Gadoth%xref no email%1,        Yael%Paran%xref no email%1]"
$ws.Range("E23").Value = "[Claudio%Liguori%NULL%3,               Mariangela%Pierantozzi%NULL%1,               Matteo%Spanetta%NULL%1,               Loredana%Sarmati%NULL%1,               Novella%Cesta%NULL%1,               Marco%Iannetta%NULL%1,               Josuel%Ora%NULL%1,               Grazia Genga%Mina%NULL%1,               Ermanno%Puxeddu%NULL%1,               Ottavia%Balbi%NULL%1,               Gabriella%Pezzuto%NULL%1,               Andrea%Magrini%NULL%1,               Paola%Rogliani%NULL%1,               Massimo%Andreoni%NULL%1,               Nicola Biagio%Mercuri%NULL%1]"
$ws.Range("E24").Value = "[Ling%Mao%xref no email%0,        Huijuan%Jin%xref no email%0,        Mengdie%Wang%xref no email%0,        Yu%Hu%xref no email%0,        Shengcai%Chen%xref no email%0,        Quanwei%He%xref no email%0,        Jiang%Chang%xref no email%0,        Candong%Hong%xref no email%0,        Yifan%Zhou%xref no email%0,        David%Wang%xref no email%0,        Xiaoping%Miao%xref no email%0,        Yanan%Li%xref no email%0,        Bo%Hu%xref no email%0]"
$ws.Range("E25").Value = "[Carmen%Martín Carreras‐Presas%carmen.martin2@universidadeuropea.es%1,               Juan%Amaro Sánchez%NULL%2,               Juan%Amaro Sánchez%NULL%0,               Antonio Francisco%López‐Sánchez%NULL%1,               Enric%Jané‐Salas%NULL%1,               Maria Luisa%Somacarrera Pérez%NULL%2,               Maria Luisa%Somacarrera Pérez%NULL%0]"
$ws.Range("E26").Value = "[Simone%Meini%simonemeini2@gmail.com%1,               Lorenzo Roberto%Suardi%NULL%2,               Lorenzo Roberto%Suardi%NULL%0,               Michele%Busoni%NULL%1,               Anna Teresa%Roberts%NULL%1,               Alberto%Fortini%NULL%1]"
$ws.Range("E28").Value = "[Muayad A.%Merza%NULL%1,               Azad A.%Haleem Al Mezori%NULL%1,               Hakar Mustafa%Mohammed%NULL%1,               Deldar Morad%Abdulah%NULL%1]"
$ws.Range("E29").Value = "[Ji Yun%Noh%NULL%2,               Jin Gu%Yoon%NULL%3,               Hye%Seong%NULL%3,               Won Suk%Choi%NULL%1,               Jang Wook%Sohn%NULL%1,               Hee Jin%Cheong%NULL%4,               Woo Joo%Kim%NULL%4,               Joon Young%Song%NULL%4]"
$ws.Range("E30").Value = "[Alberto%Paderno%NULL%1,               Alberto%Schreiber%NULL%1,               Alberto%Grammatica%albertogrammatica@libero.it%1,               Elena%Raffetti%NULL%2,               Elena%Raffetti%NULL%0,               Michele%Tomasoni%NULL%1,               Tommaso%Gualtieri%NULL%1,               Stefano%Taboni%NULL%1,               Silvia%Zorzi%NULL%1,               Davide%Lombardi%NULL%1,               Alberto%Deganello%NULL%1,               Luca Oscar%Redaelli De Zinis%NULL%1,               Roberto%Maroldi%NULL%1,               Davide%Mattavelli%NULL%1]"
$ws.Range("E31").Value = "[A.%Patel%NULL%1,               E.%Charani%NULL%1,               D.%Ariyanayagam%NULL%1,               A.%Abdulaal%NULL%1,               S.J.%Denny%NULL%1,               N.%Mughal%NULL%1,               L.S.P.%Moore%NULL%1]"
$ws.Range("E32").Value = "[Bayushi Eka%Putra%NULL%1,               Suko%Adiarto%NULL%1,               Santi Rahayu%Dewayanti%NULL%1,               Dafsah Arifa%Juzar%NULL%1]"
$ws.Range("E33").Value = "[Carlos Manuel%Romero-S\u00e1nchez%xref no email%1,        Inmaculada%D\u00edaz-Maroto%xref no email%1,        Eva%Fern\u00e1ndez-D\u00edaz%xref no email%1,        \u00c1lvaro%S\u00e1nchez-Larsen%xref no email%1,        Almudena%Layos-Romero%xref no email%1,        Jorge%Garc\u00eda-Garc\u00eda%xref no email%1,        Esther%Gonz\u00e1lez%xref no email%1,        Inmaculada%Redondo-Pe\u00f1as%xref no email%1,        Ana Bel\u00e9n%Perona-Moratalla%xref no email%1,        Jos\u00e9 Antonio%Del Valle-P\u00e9rez%xref no email%1,        Julia%Gracia-Gil%xref no email%1,        Laura%Rojas-Bartolom\u00e9%xref no email%1,        Inmaculada%Feria-Vilar%xref no email%1,        Mar\u00eda%Monteagudo%xref no email%1,        Mar\u00eda%Palao%xref no email%1,        Elena%Palaz\u00f3n-Garc\u00eda%xref no email%1,        Cristian%Alcahut-Rodr\u00edguez%xref no email%1,        David%Sopelana-Garay%xref no email%1,        Y\u00f3scar%Moreno%xref no email%1,        Javaad%Ahmad%xref no email%1,        Tom\u00e1s%Segura%xref no email%1]"
$ws.Range("E34").Value = "[İbrahim%Sayin%NULL%1,               Kadriye Kart%Yaşar%NULL%1,               Zahide Mine%Yazici%NULL%1]"
$ws.Range("E35").Value = "[Rados\u0142aw%Sierpi\u0144ski%xref no email%1,        Jaros\u0142aw%Pinkas%xref no email%1,        Mateusz%Jankowski%xref no email%1,        Wojciech S.%Zgliczy\u0144ski%xref no email%1,        Waldemar%Wierzba%xref no email%1,        Mariusz%Gujski%xref no email%1,        \u0141ukasz%Szumowski%xref no email%1]"
$ws.Range("E36").Value = "[Ciro Dantas%Soares%NULL%1,               Rejane Andrade%de Carvalho%NULL%1,               Kalline Andrade%de Carvalho%NULL%1,               Maria Goretti Freire%de Carvalho%NULL%1,               Oslei Paes%de Almeida%NULL%1]"
$ws.Range("E37").Value = "[Marlene M.%Speth%NULL%1,               Thirza%Singer-Cornelius%NULL%1,               Michael%Oberle%NULL%1,               Isabelle%Gengler%NULL%1,               Steffi J.%Brockmeier%NULL%1,               Ahmad R.%Sedaghat%NULL%1]"
$ws.Range("E38").Value = "[Ibrahim%Sultan%sultani@upmc.edu%1,               Andreas%Habertheuer%NULL%2,               Andreas%Habertheuer%NULL%0,               Asad A.%Usman%NULL%2,               Asad A.%Usman%NULL%0,               Arman%Kilic%NULL%1,               Eric%Gnall%NULL%2,               Eric%Gnall%NULL%0,               Michael E.%Friscia%NULL%1,               Dmitriy%Zubkus%NULL%1,               Hitoshi%Hirose%NULL%2,               Pablo%Sanchez%NULL%1,               Olugbenga%Okusanya%NULL%1,               Wilson Y.%Szeto%NULL%1,               Jacob%Gutsche%NULL%1]"
$ws.Range("E39").Value = "[Luigi Angelo%Vaira%luigi.vaira@gmail.com%2,               Claire%Hopkins%NULL%2,               Claire%Hopkins%NULL%0,               Giovanni%Salzano%NULL%2,               Marzia%Petrocelli%NULL%2,               Andrea%Melis%NULL%1,               Marco%Cucurullo%NULL%1,               Mario%Ferrari%NULL%1,               Laura%Gagliardini%NULL%1,               Carlotta%Pipolo%NULL%1,               Giovanna%Deiana%NULL%3,               Vito%Fiore%NULL%1,               Andrea%De Vito%NULL%0,               Nicola%Turra%NULL%1,               Sara%Canu%NULL%1,               Angelantonio%Maglio%NULL%1,               Antonello%Serra%NULL%3,               Antonello%Serra%NULL%0,               Francesco%Bussu%NULL%2,               Giordano%Madeddu%NULL%0,               Sergio%Babudieri%NULL%0,               Alessandro%Giuseppe Fois%NULL%1,               Pietro%Pirina%NULL%0,               Francesco A.%Salzano%NULL%1,               Pierluigi%De Riu%NULL%1,               Federico%Biglioli%NULL%1,               Giacomo%De Riu%NULL%2]"
$ws.Range("E40").Value = "[Carol H.%Yan%c1yan@health.ucsd.edu%3,               Farhoud%Faraji%NULL%6,               Farhoud%Faraji%NULL%0,               Divya P.%Prajapati%NULL%4,               Benjamin T.%Ostrander%NULL%2,               Adam S.%DeConde%NULL%4]"
$ws.Range("E41").Value = "[Souheil%Zayet%souhail.zayet@gmail.com%1,               Timothee%Klopfenstein%NULL%1,               Julien%Mercier%NULL%1,               N’dri Juliette%Kadiane-Oussou%NULL%1,               Ludovic%Lan Cheong Wah%NULL%1,               Pierre-Yves%Royer%NULL%1,               Lynda%Toko%NULL%1,               Vincent%Gendrin%NULL%1]"
